$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D7").Value = -7.100099999999998
$ws.Range("A8").Value = -22.48370000000001
$ws.Range("A10").Value = -21.8796
$ws.Range("A12").Value = -21.54220000000001
$ws.Range("D15").Value = -7.959599999999998
$ws.Range("A18").Value = -22.40220000000002
$ws.Range("D18").Value = -8.010599999999995
$ws.Range("D20").Value = -7.841799999999997
$ws.Range("D29").Value = -7.159899999999999
$ws.Range("D30").Value = -7.259500000000007
$ws.Range("D31").Value = -8.484500000000002
$ws.Range("A37").Value = -19.6591
$ws.Range("D40").Value = -7.988599999999994
$ws.Range("D50").Value = -8.061699999999997
$ws.Range("A55").Value = -22.30630000000001
$ws.Range("A68").Value = -21.5372
$ws.Range("D68").Value = -6.988599999999995
$ws.Range("D76").Value = -7.268799999999998
$ws.Range("A77").Value = -21.04099999999999
$ws.Range("A78").Value = -20.87059999999998
$ws.Range("A81").Value = -21.8101
$ws.Range("A82").Value = -22.1387
$ws.Range("D87").Value = -8.004199999999997
$ws.Range("D88").Value = -7.391699999999995
$ws.Range("D96").Value = -7.495000000000005
$ws.Range("D98").Value = -8.423900000000007
$ws.Range("D101").Value = -7.866499999999999
$ws.Range("D102").Value = -7.802799999999996
